$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 52: convert B52 and G52:AV52 from text to numeric values ---
$ws.Range("B52").Value = 34
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 3
$ws.Range("I52").Value = 2
$ws.Range("J52").Value = 2
$ws.Range("K52").Value = 3
$ws.Range("L52").Value = 2
$ws.Range("M52").Value = 2
$ws.Range("N52").Value = 2
$ws.Range("O52").Value = 1
$ws.Range("P52").Value = 2
$ws.Range("Q52").Value = 3
$ws.Range("R52").Value = 2
$ws.Range("S52").Value = 2
$ws.Range("T52").Value = 3
$ws.Range("U52").Value = 2
$ws.Range("V52").Value = 1
$ws.Range("W52").Value = 2
$ws.Range("X52").Value = 3
$ws.Range("Y52").Value = 2
$ws.Range("Z52").Value = 3
$ws.Range("AA52").Value = 2
$ws.Range("AB52").Value = 2
$ws.Range("AC52").Value = 3
$ws.Range("AD52").Value = 2
$ws.Range("AE52").Value = 1
$ws.Range("AF52").Value = 2
$ws.Range("AG52").Value = 3
$ws.Range("AH52").Value = 2
$ws.Range("AI52").Value = 2
$ws.Range("AJ52").Value = 2
$ws.Range("AK52").Value = 3
$ws.Range("AL52").Value = 2
$ws.Range("AM52").Value = 2
$ws.Range("AN52").Value = 2
$ws.Range("AO52").Value = 1
$ws.Range("AP52").Value = 2
$ws.Range("AQ52").Value = 2
$ws.Range("AR52").Value = 3
$ws.Range("AS52").Value = 2
$ws.Range("AT52").Value = 1
$ws.Range("AU52").Value = 2
$ws.Range("AV52").Value = 2

# --- Rows 53-58: new rows, numeric for B,G:AV; text for A,C,D,E,F ---
# Row 53
$ws.Range("A53").Value = "2025-05-24 22:59:53"
$ws.Range("B53").Value = 34
$ws.Range("C53").Value = "Bali, Indonesia"
$ws.Range("D53").Value = "SMA/SMK"
$ws.Range("E53").Value = "male"
$ws.Range("F53").Value = "jimmnnkk"
$ws.Range("G53").Value = 3
$ws.Range("H53").Value = 2
$ws.Range("I53").Value = 3
$ws.Range("J53").Value = 2
$ws.Range("K53").Value = 2
$ws.Range("L53").Value = 3
$ws.Range("M53").Value = 2
$ws.Range("N53").Value = 3
$ws.Range("O53").Value = 2
$ws.Range("P53").Value = 2
$ws.Range("Q53").Value = 1
$ws.Range("R53").Value = 2
$ws.Range("S53").Value = 2
$ws.Range("T53").Value = 1
$ws.Range("U53").Value = 2
$ws.Range("V53").Value = 2
$ws.Range("W53").Value = 1
$ws.Range("X53").Value = 2
$ws.Range("Y53").Value = 2
$ws.Range("Z53").Value = 1
$ws.Range("AA53").Value = 2
$ws.Range("AB53").Value = 3
$ws.Range("AC53").Value = 2
$ws.Range("AD53").Value = 2
$ws.Range("AE53").Value = 3
$ws.Range("AF53").Value = 2
$ws.Range("AG53").Value = 2
$ws.Range("AH53").Value = 3
$ws.Range("AI53").Value = 2
$ws.Range("AJ53").Value = 2
$ws.Range("AK53").Value = 3
$ws.Range("AL53").Value = 2
$ws.Range("AM53").Value = 2
$ws.Range("AN53").Value = 3
$ws.Range("AO53").Value = 2
$ws.Range("AP53").Value = 2
$ws.Range("AQ53").Value = 1
$ws.Range("AR53").Value = 2
$ws.Range("AS53").Value = 3
$ws.Range("AT53").Value = 2
$ws.Range("AU53").Value = 2
$ws.Range("AV53").Value = 1

# Row 54
$ws.Range("A54").Value = "2025-05-26 11:37:50"
$ws.Range("B54").Value = 45
$ws.Range("C54").Value = "Jakarta"
$ws.Range("D54").Value = "SMA/SMK"
$ws.Range("E54").Value = "male"
$ws.Range("F54").Value = "M. Rakhmat Dramaga"
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 3
$ws.Range("I54").Value = 2
$ws.Range("J54").Value = 3
$ws.Range("K54").Value = 2
$ws.Range("L54").Value = 3
$ws.Range("M54").Value = 3
$ws.Range("N54").Value = 3
$ws.Range("O54").Value = 3
$ws.Range("P54").Value = 2
$ws.Range("Q54").Value = 3
$ws.Range("R54").Value = 2
$ws.Range("S54").Value = 4
$ws.Range("T54").Value = 3
$ws.Range("U54").Value = 3
$ws.Range("V54").Value = 2
$ws.Range("W54").Value = 3
$ws.Range("X54").Value = 2
$ws.Range("Y54").Value = 4
$ws.Range("Z54").Value = 3
$ws.Range("AA54").Value = 4
$ws.Range("AB54").Value = 3
$ws.Range("AC54").Value = 3
$ws.Range("AD54").Value = 2
$ws.Range("AE54").Value = 3
$ws.Range("AF54").Value = 2
$ws.Range("AG54").Value = 3
$ws.Range("AH54").Value = 2
$ws.Range("AI54").Value = 3
$ws.Range("AJ54").Value = 4
$ws.Range("AK54").Value = 3
$ws.Range("AL54").Value = 2
$ws.Range("AM54").Value = 3
$ws.Range("AN54").Value = 2
$ws.Range("AO54").Value = 3
$ws.Range("AP54").Value = 3
$ws.Range("AQ54").Value = 2
$ws.Range("AR54").Value = 3
$ws.Range("AS54").Value = 3
$ws.Range("AT54").Value = 2
$ws.Range("AU54").Value = 2
$ws.Range("AV54").Value = 3

# Row 55
$ws.Range("A55").Value = "2025-05-26 11:48:38"
$ws.Range("B55").Value = 45
$ws.Range("C55").Value = "Jakarta"
$ws.Range("D55").Value = "SMA/SMK"
$ws.Range("E55").Value = "male"
$ws.Range("F55").Value = "M. Rakhmat Dramaga"
$ws.Range("G55").Value = 4
$ws.Range("H55").Value = 4
$ws.Range("I55").Value = 4
$ws.Range("J55").Value = 4
$ws.Range("K55").Value = 4
$ws.Range("L55").Value = 4
$ws.Range("M55").Value = 4
$ws.Range("N55").Value = 4
$ws.Range("O55").Value = 4
$ws.Range("P55").Value = 4
$ws.Range("Q55").Value = 4
$ws.Range("R55").Value = 4
$ws.Range("S55").Value = 4
$ws.Range("T55").Value = 4
$ws.Range("U55").Value = 4
$ws.Range("V55").Value = 4
$ws.Range("W55").Value = 4
$ws.Range("X55").Value = 4
$ws.Range("Y55").Value = 4
$ws.Range("Z55").Value = 4
$ws.Range("AA55").Value = 4
$ws.Range("AB55").Value = 4
$ws.Range("AC55").Value = 4
$ws.Range("AD55").Value = 4
$ws.Range("AE55").Value = 4
$ws.Range("AF55").Value = 4
$ws.Range("AG55").Value = 4
$ws.Range("AH55").Value = 4
$ws.Range("AI55").Value = 4
$ws.Range("AJ55").Value = 4
$ws.Range("AK55").Value = 4
$ws.Range("AL55").Value = 4
$ws.Range("AM55").Value = 4
$ws.Range("AN55").Value = 4
$ws.Range("AO55").Value = 4
$ws.Range("AP55").Value = 4
$ws.Range("AQ55").Value = 4
$ws.Range("AR55").Value = 4
$ws.Range("AS55").Value = 4
$ws.Range("AT55").Value = 4
$ws.Range("AU55").Value = 4
$ws.Range("AV55").Value = 4

# Row 56
$ws.Range("A56").Value = "2025-05-26 11:57:52"
$ws.Range("B56").Value = 23
$ws.Range("C56").Value = "Medan"
$ws.Range("D56").Value = "SMA/SMK"
$ws.Range("E56").Value = "male"
$ws.Range("F56").Value = "Bahteramon"
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 3
$ws.Range("I56").Value = 3
$ws.Range("J56").Value = 2
$ws.Range("K56").Value = 3
$ws.Range("L56").Value = 2
$ws.Range("M56").Value = 3
$ws.Range("N56").Value = 3
$ws.Range("O56").Value = 2
$ws.Range("P56").Value = 3
$ws.Range("Q56").Value = 3
$ws.Range("R56").Value = 2
$ws.Range("S56").Value = 3
$ws.Range("T56").Value = 4
$ws.Range("U56").Value = 3
$ws.Range("V56").Value = 4
$ws.Range("W56").Value = 3
$ws.Range("X56").Value = 3
$ws.Range("Y56").Value = 4
$ws.Range("Z56").Value = 3
$ws.Range("AA56").Value = 3
$ws.Range("AB56").Value = 2
$ws.Range("AC56").Value = 3
$ws.Range("AD56").Value = 4
$ws.Range("AE56").Value = 3
$ws.Range("AF56").Value = 2
$ws.Range("AG56").Value = 3
$ws.Range("AH56").Value = 4
$ws.Range("AI56").Value = 3
$ws.Range("AJ56").Value = 3
$ws.Range("AK56").Value = 2
$ws.Range("AL56").Value = 3
$ws.Range("AM56").Value = 3
$ws.Range("AN56").Value = 4
$ws.Range("AO56").Value = 3
$ws.Range("AP56").Value = 2
$ws.Range("AQ56").Value = 3
$ws.Range("AR56").Value = 4
$ws.Range("AS56").Value = 3
$ws.Range("AT56").Value = 2
$ws.Range("AU56").Value = 3
$ws.Range("AV56").Value = 3

# Row 57
$ws.Range("A57").Value = "2025-05-26 15:01:19"
$ws.Range("B57").Value = 23
$ws.Range("C57").Value = "Medan"
$ws.Range("D57").Value = "SMA/SMK"
$ws.Range("E57").Value = "male"
$ws.Range("F57").Value = "Bahteramon"
$ws.Range("G57").Value = 3
$ws.Range("H57").Value = 3
$ws.Range("I57").Value = 3
$ws.Range("J57").Value = 2
$ws.Range("K57").Value = 3
$ws.Range("L57").Value = 3
$ws.Range("M57").Value = 3
$ws.Range("N57").Value = 2
$ws.Range("O57").Value = 3
$ws.Range("P57").Value = 2
$ws.Range("Q57").Value = 3
$ws.Range("R57").Value = 3
$ws.Range("S57").Value = 3
$ws.Range("T57").Value = 3
$ws.Range("U57").Value = 3
$ws.Range("V57").Value = 2
$ws.Range("W57").Value = 3
$ws.Range("X57").Value = 2
$ws.Range("Y57").Value = 3
$ws.Range("Z57").Value = 2
$ws.Range("AA57").Value = 3
$ws.Range("AB57").Value = 2
$ws.Range("AC57").Value = 3
$ws.Range("AD57").Value = 2
$ws.Range("AE57").Value = 3
$ws.Range("AF57").Value = 2
$ws.Range("AG57").Value = 3
$ws.Range("AH57").Value = 2
$ws.Range("AI57").Value = 3
$ws.Range("AJ57").Value = 2
$ws.Range("AK57").Value = 3
$ws.Range("AL57").Value = 2
$ws.Range("AM57").Value = 3
$ws.Range("AN57").Value = 3
$ws.Range("AO57").Value = 3
$ws.Range("AP57").Value = 3
$ws.Range("AQ57").Value = 4
$ws.Range("AR57").Value = 3
$ws.Range("AS57").Value = 2
$ws.Range("AT57").Value = 3
$ws.Range("AU57").Value = 2
$ws.Range("AV57").Value = 3

# Row 58
$ws.Range("A58").Value = "2025-05-26 15:12:41"
$ws.Range("B58").Value = 23
$ws.Range("C58").Value = "Medan"
$ws.Range("D58").Value = "SMA/SMK"
$ws.Range("E58").Value = "male"
$ws.Range("F58").Value = "Bahteramon"
$ws.Range("G58").Value = 3
$ws.Range("H58").Value = 3
$ws.Range("I58").Value = 2
$ws.Range("J58").Value = 3
$ws.Range("K58").Value = 3
$ws.Range("L58").Value = 2
$ws.Range("M58").Value = 3
$ws.Range("N58").Value = 2
$ws.Range("O58").Value = 3
$ws.Range("P58").Value = 2
$ws.Range("Q58").Value = 3
$ws.Range("R58").Value = 4
$ws.Range("S58").Value = 3
$ws.Range("T58").Value = 2
$ws.Range("U58").Value = 3
$ws.Range("V58").Value = 2
$ws.Range("W58").Value = 3
$ws.Range("X58").Value = 2
$ws.Range("Y58").Value = 3
$ws.Range("Z58").Value = 2
$ws.Range("AA58").Value = 3
$ws.Range("AB58").Value = 2
$ws.Range("AC58").Value = 3
$ws.Range("AD58").Value = 4
$ws.Range("AE58").Value = 3
$ws.Range("AF58").Value = 2
$ws.Range("AG58").Value = 3
$ws.Range("AH58").Value = 2
$ws.Range("AI58").Value = 3
$ws.Range("AJ58").Value = 4
$ws.Range("AK58").Value = 3
$ws.Range("AL58").Value = 2
$ws.Range("AM58").Value = 3
$ws.Range("AN58").Value = 2
$ws.Range("AO58").Value = 3
$ws.Range("AP58").Value = 2
$ws.Range("AQ58").Value = 3
$ws.Range("AR58").Value = 2
$ws.Range("AS58").Value = 3
$ws.Range("AT58").Value = 2
$ws.Range("AU58").Value = 3
$ws.Range("AV58").Value = 2

# --- Row 59: new row, all text (inline string type) ---
# A,C,D,E,F are naturally non-numeric text. B and G:AV are numeric-looking
# strings, so a leading apostrophe is used to force text storage, and the
# resulting quote-prefix style is then cleared back to Normal so no stray
# style is left on the cells (matching the source which has no "s" attr).
$ws.Range("A59").Value = "2025-05-26 15:20:11"
$ws.Range("B59").Value = "'23"
$ws.Range("C59").Value = "Medan"
$ws.Range("D59").Value = "SMP"
$ws.Range("E59").Value = "male"
$ws.Range("F59").Value = "Bahteramon"
$ws.Range("G59").Value = "'3"
$ws.Range("H59").Value = "'4"
$ws.Range("I59").Value = "'3"
$ws.Range("J59").Value = "'4"
$ws.Range("K59").Value = "'3"
$ws.Range("L59").Value = "'3"
$ws.Range("M59").Value = "'2"
$ws.Range("N59").Value = "'3"
$ws.Range("O59").Value = "'2"
$ws.Range("P59").Value = "'3"
$ws.Range("Q59").Value = "'4"
$ws.Range("R59").Value = "'3"
$ws.Range("S59").Value = "'2"
$ws.Range("T59").Value = "'3"
$ws.Range("U59").Value = "'3"
$ws.Range("V59").Value = "'2"
$ws.Range("W59").Value = "'3"
$ws.Range("X59").Value = "'2"
$ws.Range("Y59").Value = "'3"
$ws.Range("Z59").Value = "'2"
$ws.Range("AA59").Value = "'3"
$ws.Range("AB59").Value = "'2"
$ws.Range("AC59").Value = "'3"
$ws.Range("AD59").Value = "'2"
$ws.Range("AE59").Value = "'3"
$ws.Range("AF59").Value = "'2"
$ws.Range("AG59").Value = "'3"
$ws.Range("AH59").Value = "'2"
$ws.Range("AI59").Value = "'3"
$ws.Range("AJ59").Value = "'2"
$ws.Range("AK59").Value = "'3"
$ws.Range("AL59").Value = "'2"
$ws.Range("AM59").Value = "'2"
$ws.Range("AN59").Value = "'2"
$ws.Range("AO59").Value = "'2"
$ws.Range("AP59").Value = "'2"
$ws.Range("AQ59").Value = "'2"
$ws.Range("AR59").Value = "'2"
$ws.Range("AS59").Value = "'2"
$ws.Range("AT59").Value = "'2"
$ws.Range("AU59").Value = "'2"
$ws.Range("AV59").Value = "'2"
$ws.Range("A59:AV59").Style = "Normal"
